# Trade #84 closed at 2026-02-17 15:53:15 - unknown UNKNOWN +0.000%
#
# Applies the update for the newly closed MarketMaking trade #84:
#   - Appends the trade row to "All Trades" and "MarketMaking" sheets
#   - Refreshes the aggregate numbers on "Summary" and "Strategy Status"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.07   # Current Capital
$summary.Range("B4").Value = 0.06      # Total P&L $
$summary.Range("B5").Value = 0.01      # Total P&L %
$summary.Range("B6").Value = 84        # Total Trades
$summary.Range("B7").Value = 29        # Winning Trades
$summary.Range("B9").Value = 34.52     # Win Rate %

# ---------------------------------------------------------------------
# Sheet: Strategy Status (MarketMaking row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.07     # Capital
$status.Range("D4").Value = 84         # Trades
$status.Range("E4").Value = 0.06       # P&L $
$status.Range("F4").Value = 0.07       # P&L %
$status.Range("G4").Value = 34.52      # Win Rate %

# ---------------------------------------------------------------------
# New closed-trade row shared by "All Trades" and "MarketMaking"
# ---------------------------------------------------------------------
$tradeSheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $tradeSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 85

    $ws.Cells.Item($row, 1).Value = 84                # Trade #

    # Date / Time columns must stay literal text, not get auto-converted
    # to date/time serial numbers.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"       # Date
    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = "15:53:09"         # Time

    $ws.Cells.Item($row, 4).Value = "MarketMaking"     # Strategy
    $ws.Cells.Item($row, 5).Value = "DOWN"             # Side
    $ws.Cells.Item($row, 6).Value = 0.18               # Entry Price
    $ws.Cells.Item($row, 7).Value = 0.26               # Exit Price
    $ws.Cells.Item($row, 8).Value = "CLOSED"           # Status
    $ws.Cells.Item($row, 9).Value = 44.4444            # P&L %
    $ws.Cells.Item($row, 10).Value = 0.08              # P&L $
    $ws.Cells.Item($row, 11).Value = 100.07            # Capital After
    $ws.Cells.Item($row, 12).Value = 0                 # Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                 # Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6               # Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item($row, 16).Value = "early_exit"      # Exit Reason
    $ws.Cells.Item($row, 17).Value = 0.14              # Duration (min)
}
